# Update "想去人数" (want-to-go count) values in column F across the
# workbook's sheets, matching the data refresh recorded in the diff.

$wb = $excel.ActiveWorkbook

# Worksheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 1151
$ws1.Range("F5").Value  = 10
$ws1.Range("F6").Value  = 2773
$ws1.Range("F7").Value  = 231
$ws1.Range("F9").Value  = 97
$ws1.Range("F10").Value = 287
$ws1.Range("F11").Value = 199
$ws1.Range("F12").Value = 708
$ws1.Range("F13").Value = 108
$ws1.Range("F14").Value = 134
$ws1.Range("F15").Value = 1721
$ws1.Range("F16").Value = 311
$ws1.Range("F17").Value = 46

# Worksheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value  = 23
$ws2.Range("F10").Value = 37
$ws2.Range("F11").Value = 44
$ws2.Range("F18").Value = 36

# Worksheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6359
$ws3.Range("F3").Value = 798
$ws3.Range("F4").Value = 2031
$ws3.Range("F5").Value = 266

# Worksheet "全部类型" (All types) - combined view of the three sheets above
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 6359
$ws4.Range("F3").Value  = 798
$ws4.Range("F4").Value  = 2031
$ws4.Range("F5").Value  = 266
$ws4.Range("F12").Value = 1151
$ws4.Range("F13").Value = 10
$ws4.Range("F15").Value = 23
$ws4.Range("F17").Value = 2773
$ws4.Range("F19").Value = 231
$ws4.Range("F20").Value = 37
$ws4.Range("F21").Value = 44
$ws4.Range("F24").Value = 97
$ws4.Range("F25").Value = 287
$ws4.Range("F27").Value = 199
$ws4.Range("F28").Value = 708
$ws4.Range("F29").Value = 108
$ws4.Range("F30").Value = 134
$ws4.Range("F32").Value = 1721
$ws4.Range("F33").Value = 311
$ws4.Range("F36").Value = 46
$ws4.Range("F39").Value = 36
